# Apply the commit: "regenerate synthetic data with dcterms:description and dcterms:title properties"

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. RightsStatement sheet: update E2 (note) text
# -----------------------------------------------------------------------
$wsRights = $wb.Worksheets.Item("RightsStatement")
$wsRights.Range("E2").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."

# -----------------------------------------------------------------------
# 2. Property sheet: header rename, value normalization, and insertion of
#    two new rows (dcterms:description, dcterms:title)
# -----------------------------------------------------------------------
$wsProp = $wb.Worksheets.Item("Property")

# Header column B renamed from "cms:propertyFilterable" to "filterable"
$wsProp.Range("B1").Value = "filterable"

# Insert a new row for "dcterms:description" right before the existing
# "dcterms:extent" row (currently row 3), shifting everything else down.
$wsProp.Rows.Item(3).Insert()

# Insert a new row for "dcterms:title" right before the existing
# "dcterms:type" row. Before this second insert, "dcterms:type" sits at
# row 13 (it was row 12, plus 1 from the insert above).
$wsProp.Rows.Item(13).Insert()

# Rewrite the B column ("filterable") for every existing data row from
# the old python repr ("<class 'filter'>") to a plain boolean string, and
# fill in the two newly inserted rows.
$propertyRows = @(
    @{ Row = 2;  Id = "vra:culturalContext"; Filterable = "true";  Label = "Cultural context"; Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Cultural%20context" },
    @{ Row = 3;  Id = "dcterms:description"; Filterable = "false"; Label = "Description";       Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Description" },
    @{ Row = 4;  Id = "dcterms:extent";      Filterable = "true";  Label = "Extent";             Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Extent" },
    @{ Row = 5;  Id = "dcterms:language";    Filterable = "true";  Label = "Language";           Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Language" },
    @{ Row = 6;  Id = "vra:material";        Filterable = "true";  Label = "Material";           Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Material" },
    @{ Row = 7;  Id = "dcterms:medium";      Filterable = "true";  Label = "Medium";             Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Medium" },
    @{ Row = 8;  Id = "dcterms:publisher";   Filterable = "true";  Label = "Publisher";          Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Publisher" },
    @{ Row = 9;  Id = "dcterms:source";      Filterable = "true";  Label = "Source";             Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Source" },
    @{ Row = 10; Id = "dcterms:spatial";     Filterable = "true";  Label = "Spatial";            Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Spatial" },
    @{ Row = 11; Id = "dcterms:subject";     Filterable = "true";  Label = "Subject";            Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Subject" },
    @{ Row = 12; Id = "vra:technique";       Filterable = "true";  Label = "Technique";          Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Technique" },
    @{ Row = 13; Id = "dcterms:title";       Filterable = "false"; Label = "Title";              Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Title" },
    @{ Row = 14; Id = "dcterms:type";        Filterable = "true";  Label = "Type";               Range = "urn:paradicms_etl:pipeline:synthetic_data:property_range:Type" }
)

foreach ($row in $propertyRows) {
    $r = $row.Row
    $wsProp.Cells.Item($r, 1).Value = $row.Id
    # Prefix boolean-looking text with an apostrophe so Excel stores it as
    # plain text ("true"/"false") rather than auto-converting to a native
    # Boolean cell type.
    $wsProp.Cells.Item($r, 2).Value = "'" + $row.Filterable
    $wsProp.Cells.Item($r, 3).Value = "urn:paradicms_etl:pipeline:synthetic_data:property_group"
    $wsProp.Cells.Item($r, 4).Value = $row.Label
    $wsProp.Cells.Item($r, 5).Value = $row.Range
    $wsProp.Cells.Item($r, 6).Value = "'true"
}

# -----------------------------------------------------------------------
# 3. Person sheet: shuffle the "page" (F column) links between Wikipedia
#    and Wikidata for rows 3, 4, 5, 6 (row 2 is unchanged).
# -----------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("F3").Value = "http://www.wikidata.org/entity/Q7251"
$wsPerson.Range("F4").Value = "http://www.wikidata.org/entity/Q7251"
$wsPerson.Range("F5").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$wsPerson.Range("F6").Value = "http://www.wikidata.org/entity/Q7251"
